$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '34.475.32'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.03%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.802.45'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.40%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.41%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '224.16'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -1.82%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.604'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.10%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.33%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '39.00'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +7.33%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.289'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -3.55%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0669'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -4.57%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0999'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +3.21%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '2.057.37'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.73%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.800.03'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -0.80%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '10.85'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -4.56%  '
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = 'Polygon'
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.629'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -3.49%  '
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = 'WrappedBTC'
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '34.438.82'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.20%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '4.36'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -2.62%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '67.93'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -2.92%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '239.88'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -2.92%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0766'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -3.48%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.05'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -4.60%  '
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.32%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.08'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -2.94%  '
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -3.39%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '170.61'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.63%  '
$ws.Range("B26").NumberFormat = "@"
$ws.Range("B26").Value = 'Cosmos'
$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.66'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -5.85%  '
$ws.Range("B27").NumberFormat = "@"
$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.47'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -1.26%  '
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -1.13%  '
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.32%  '
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -2.23%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.73'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -2.88%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0513'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -2.83%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.81'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -4.73%  '
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -0.63%  '
$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = 'ImmutableX'
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.639'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -3.72%  '
$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = 'Maker'
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.307.39'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -6.32%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -0.81%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0186'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -1.77%  '
$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.31'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -3.63%  '
$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = 'WEMIXToken'
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.23'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +2.93%  '
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = 'HuobiToken'
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.44'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +0.58%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '81.91'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -1.03%  '
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -1.95%  '
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = 'ARBITRUM'
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.941'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -2.29%  '
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = 'InjectiveProtocol'
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '14.13'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +3.76%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0517'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +3.75%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.960.48'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -0.64%  '
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -0.37%  '
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -5.39%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '101.54'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -1.83%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0611'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -1.06%  '
